# Auto-generated edit script: refresh cryptos price table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.028.35"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.866.17"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.68%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "599.16"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.81%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "170.39"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.859.36"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.97%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.529"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.41"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.455"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.26%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000259"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.92"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.518.99"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.868.77"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.177.45"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.09"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.12%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.32"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.20%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.76"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "465.96"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -6.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.737"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.68%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000158"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -5.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.07"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.50%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.47%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.15%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.999"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.96"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.96"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.020.21"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.70"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.13"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.48"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.834.18"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.62%  "
$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.104"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.87%  "
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.70"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +12.46%  "
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.140"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.02"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.52%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.96%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.311"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.44%  "
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "423.66"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -2.44%  "
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "FLOKI"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.000297"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.99%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.60"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "47.09"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.35"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "142.33"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.55%  "
